$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-21 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-22 Monday", 2)

$d.Content.Find.Execute("140÷7=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "335÷9=37, 2", 2)
$d.Content.Find.Execute("299÷5=59, 4", $true, $false, $false, $false, $false, $true, 1, $false, "833÷9=92, 5", 2)
$d.Content.Find.Execute("733÷7=104, 5", $true, $false, $false, $false, $false, $true, 1, $false, "700÷5=140, 0", 2)
$d.Content.Find.Execute("676÷3=225, 1", $true, $false, $false, $false, $false, $true, 1, $false, "724÷6=120, 4", 2)
$d.Content.Find.Execute("779÷2=389, 1", $true, $false, $false, $false, $false, $true, 1, $false, "961÷2=480, 1", 2)

$d.Content.Find.Execute("868÷8=108, 4", $true, $false, $false, $false, $false, $true, 1, $false, "808÷5=161, 3", 2)
$d.Content.Find.Execute("624÷8=78, 0", $true, $false, $false, $false, $false, $true, 1, $false, "932÷8=116, 4", 2)
$d.Content.Find.Execute("148÷8=18, 4", $true, $false, $false, $false, $false, $true, 1, $false, "933÷9=103, 6", 2)
$d.Content.Find.Execute("461÷7=65, 6", $true, $false, $false, $false, $false, $true, 1, $false, "384÷3=128, 0", 2)
$d.Content.Find.Execute("213÷7=30, 3", $true, $false, $false, $false, $false, $true, 1, $false, "759÷7=108, 3", 2)

$d.Content.Find.Execute("964÷9=107, 1", $true, $false, $false, $false, $false, $true, 1, $false, "754÷5=150, 4", 2)
$d.Content.Find.Execute("815÷3=271, 2", $true, $false, $false, $false, $false, $true, 1, $false, "131÷3=43, 2", 2)
$d.Content.Find.Execute("998÷9=110, 8", $true, $false, $false, $false, $false, $true, 1, $false, "871÷9=96, 7", 2)
$d.Content.Find.Execute("976÷3=325, 1", $true, $false, $false, $false, $false, $true, 1, $false, "981÷7=140, 1", 2)
$d.Content.Find.Execute("734÷6=122, 2", $true, $false, $false, $false, $false, $true, 1, $false, "846÷8=105, 6", 2)

$d.Content.Find.Execute("735÷6=122, 3", $true, $false, $false, $false, $false, $true, 1, $false, "647÷7=92, 3", 2)
$d.Content.Find.Execute("322÷7=46, 0", $true, $false, $false, $false, $false, $true, 1, $false, "830÷5=166, 0", 2)
$d.Content.Find.Execute("310÷5=62, 0", $true, $false, $false, $false, $false, $true, 1, $false, "629÷7=89, 6", 2)
$d.Content.Find.Execute("172÷7=24, 4", $true, $false, $false, $false, $false, $true, 1, $false, "588÷2=294, 0", 2)
$d.Content.Find.Execute("664÷6=110, 4", $true, $false, $false, $false, $false, $true, 1, $false, "690÷6=115, 0", 2)

$d.Content.Find.Execute("895÷4=223, 3", $true, $false, $false, $false, $false, $true, 1, $false, "994÷9=110, 4", 2)
$d.Content.Find.Execute("251÷4=62, 3", $true, $false, $false, $false, $false, $true, 1, $false, "912÷6=152, 0", 2)
$d.Content.Find.Execute("814÷3=271, 1", $true, $false, $false, $false, $false, $true, 1, $false, "123÷3=41, 0", 2)
$d.Content.Find.Execute("116÷9=12, 8", $true, $false, $false, $false, $false, $true, 1, $false, "282÷6=47, 0", 2)
$d.Content.Find.Execute("975÷3=325, 0", $true, $false, $false, $false, $false, $true, 1, $false, "620÷9=68, 8", 2)
